# P3518_MA_octahedron_patterning_v1.xlsx
#
# The "Names" sheet previously kept each cargo strand's internal numeric
# ID (e.g. "_id1", "_id3") embedded in its name, and used a standalone
# regex (stored in A1) to parse that naming convention back apart. Now
# that this workbook is wired into the generic plate-import system, the
# "_idN" tokens are no longer needed in the displayed names (the plate
# system tracks identity itself) and the "-10mer" size qualifier moves
# from a suffix to a prefix ahead of the "_h#_pos#" part. A1's helper
# regex is replaced by a plain "name-side-position" label describing the
# new convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Names")

# Column header / naming-convention description
$ws.Range("A1").Value = "name-side-position"

# Row 2 - antiBart (h2)
$ws.Range("B2").Value = "antiBart_h2_pos2"
$ws.Range("C2").Value = "antiBart_h2_pos3"
$ws.Range("D2").Value = "antiBart_h2_pos8"
$ws.Range("E2").Value = "antiBart_h2_pos9"
$ws.Range("F2").Value = "antiBart_h2_pos14"
$ws.Range("G2").Value = "antiBart_h2_pos15"
$ws.Range("H2").Value = "antiBart_h2_pos20"
$ws.Range("I2").Value = "antiBart_h2_pos21"
$ws.Range("J2").Value = "antiBart_h2_pos26"
$ws.Range("K2").Value = "antiBart_h2_pos27"

# Row 3 - antiEdna (h2)
$ws.Range("B3").Value = "antiEdna_h2_pos5"
$ws.Range("C3").Value = "antiEdna_h2_pos6"
$ws.Range("D3").Value = "antiEdna_h2_pos11"
$ws.Range("E3").Value = "antiEdna_h2_pos12"
$ws.Range("F3").Value = "antiEdna_h2_pos17"
$ws.Range("G3").Value = "antiEdna_h2_pos18"
$ws.Range("H3").Value = "antiEdna_h2_pos23"
$ws.Range("I3").Value = "antiEdna_h2_pos24"
$ws.Range("J3").Value = "antiEdna_h2_pos29"
$ws.Range("K3").Value = "antiEdna_h2_pos30"

# Row 4 - 10mer cargo strands (h2)
$ws.Range("B4").Value = "antiSmithers-10mer_h2_pos4"
$ws.Range("C4").Value = "antiQuimby-10mer_h2_pos7"
$ws.Range("D4").Value = "antiPatty-10mer_h2_pos10"
$ws.Range("E4").Value = "antiMarge-10mer_h2_pos13"
$ws.Range("F4").Value = "antiSmithers-10mer_h2_pos14"
$ws.Range("G4").Value = "antiLisa-10mer_h2_pos16"
$ws.Range("H4").Value = "antiQuimby-10mer_h2_pos17"
$ws.Range("I4").Value = "antiKrusty-10mer_h2_pos19"
$ws.Range("J4").Value = "antiPatty-10mer_h2_pos20"
$ws.Range("K4").Value = "antiHomer-10mer_h2_pos22"
$ws.Range("L4").Value = "antiMarge-10mer_h2_pos23"
$ws.Range("M4").Value = "antiLisa-10mer_h2_pos26"
$ws.Range("N4").Value = "antiKrusty-10mer_h2_pos29"
$ws.Range("O4").Value = "antiHomer-10mer_h2_pos32"

# Row 5 - 10mer cargo strands (h5)
$ws.Range("B5").Value = "Homer-10mer_h5_pos1"
$ws.Range("C5").Value = "Krusty-10mer_h5_pos6"
$ws.Range("D5").Value = "Lisa-10mer_h5_pos11"
$ws.Range("E5").Value = "Marge-10mer_h5_pos16"
$ws.Range("F5").Value = "Patty-10mer_h5_pos21"
$ws.Range("G5").Value = "Quimby-10mer_h5_pos26"
$ws.Range("H5").Value = "Smithers-10mer_h5_pos31"

# Row 6 - biotin anchor (unchanged name, just re-seated in the shared
# string table now that it moved next to the other plate entries)
$ws.Range("B6").Value = "biotin-anchor_id17"

# The author's selection moved from A8 to A10 while reviewing the sheet.
$ws.Activate()
$ws.Range("A10").Select()
